# Apply crypto price/volume updates scraped on Thu Nov  2 14:40:13 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.085.70'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').Value = '1.830.02'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('D5').Value = '''232.64'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').Value = '''0.615'
$ws.Range('E6').Value = '  +1.94%  '
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('D8').Value = '''42.59'
$ws.Range('E8').Value = '  +4.74%  '
$ws.Range('D9').Value = '''0.309'
$ws.Range('E9').Value = '  +5.64%  '
$ws.Range('D10').Value = '''0.0688'
$ws.Range('E10').Value = '  +2.21%  '
$ws.Range('D11').Value = '''0.100'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = '2.099.97'
$ws.Range('D13').Value = '1.851.04'
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').Value = '''11.15'
$ws.Range('E14').Value = '  +2.13%  '
$ws.Range('D15').Value = '''0.663'
$ws.Range('E15').Value = '  +4.42%  '
$ws.Range('D16').Value = '''4.68'
$ws.Range('E16').Value = '  +6.13%  '
$ws.Range('D17').Value = '35.105.28'
$ws.Range('E17').Value = '  +1.25%  '
$ws.Range('D18').Value = '''70.05'
$ws.Range('E18').Value = '  +3.09%  '
$ws.Range('D19').Value = '0.0₃0789'
$ws.Range('E19').Value = '  +2.50%  '
$ws.Range('D20').Value = '''239.77'
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('D21').Value = '''11.81'
$ws.Range('E21').Value = '  +6.23%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '''1.01'
$ws.Range('E22').Value = '  +0.56%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '''4.57'
$ws.Range('E23').Value = '  +11.23%  '
$ws.Range('E24').Value = '  +3.65%  '
$ws.Range('D25').Value = '''171.46'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').Value = '''7.78'
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('D27').Value = '''17.52'
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('E29').Value = '  +29.28%  '
$ws.Range('E30').Value = '  +0.66%  '
$ws.Range('D31').Value = '3.343.88'
$ws.Range('E31').Value = '  +37.63%  '
$ws.Range('D32').Value = '''0.0552'
$ws.Range('E32').Value = '  +7.03%  '
$ws.Range('D33').Value = '''3.90'
$ws.Range('E33').Value = '  +3.07%  '
$ws.Range('D34').Value = '''4.00'
$ws.Range('E34').Value = '  +3.82%  '
$ws.Range('D35').Value = '''1.79'
$ws.Range('E35').Value = '  -0.65%  '
$ws.Range('D36').Value = '''93.35'
$ws.Range('E36').Value = '  +9.21%  '
$ws.Range('D37').Value = '''0.678'
$ws.Range('E37').Value = '  +5.42%  '
$ws.Range('E38').Value = '  +4.73%  '
$ws.Range('D39').Value = '1.325.92'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('D40').Value = '''0.0193'
$ws.Range('E40').Value = '  +2.45%  '
$ws.Range('E41').Value = '  +2.08%  '
$ws.Range('D42').Value = '''0.991'
$ws.Range('E42').Value = '  +4.93%  '
$ws.Range('D43').Value = '''2.36'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '''14.82'
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('D45').Value = '''2.45'
$ws.Range('E45').Value = '  +1.13%  '
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('D47').Value = '''6.21'
$ws.Range('E47').Value = '  +7.93%  '
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('D49').Value = '2.008.85'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('E50').Value = '  +0.59%  '
$ws.Range('D51').Value = '''100.85'
$ws.Range('E51').Value = '  -0.33%  '
